$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 6192
$ws.Range("L3").Value = 6713
$ws.Range("L4").Value = 1663
$ws.Range("L5").Value = 398
$ws.Range("L6").Value = 5529
$ws.Range("L7").Value = 20495

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L3").Value = 57
$ws.Range("L4").Value = 25
$ws.Range("L7").Value = 227

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 408
$ws.Range("L7").Value = 1350

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L2").Value = 251
$ws.Range("L3").Value = 325
$ws.Range("L5").Value = 23
$ws.Range("L7").Value = 925

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 238
$ws.Range("L7").Value = 788

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L4").Value = 23
$ws.Range("L6").Value = 97
$ws.Range("L7").Value = 401

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L3").Value = 145
$ws.Range("L6").Value = 75
$ws.Range("L7").Value = 360

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("L3").Value = 21
$ws.Range("L7").Value = 91

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L6").Value = 163
$ws.Range("L7").Value = 658
$ws.Range("L8").Value = 1350
$ws.Range("L10").Value = 136
$ws.Range("L11").Value = 338
$ws.Range("L14").Value = 101
$ws.Range("L18").Value = 140
$ws.Range("L19").Value = 554
$ws.Range("L20").Value = 522
$ws.Range("L29").Value = 1146
$ws.Range("L30").Value = 91
$ws.Range("L31").Value = 201
$ws.Range("L33").Value = 925
$ws.Range("L37").Value = 788
$ws.Range("L40").Value = 54
$ws.Range("L41").Value = 89
$ws.Range("L42").Value = 654
$ws.Range("L43").Value = 154
$ws.Range("L48").Value = 269
$ws.Range("L49").Value = 111
$ws.Range("L52").Value = 436
$ws.Range("L53").Value = 227
$ws.Range("L56").Value = 20
$ws.Range("L60").Value = 134
$ws.Range("L63").Value = 63
$ws.Range("L65").Value = 401
$ws.Range("L66").Value = 59
$ws.Range("L76").Value = 317
$ws.Range("L78").Value = 268
$ws.Range("L80").Value = 67
$ws.Range("L84").Value = 197
$ws.Range("L85").Value = 1019
$ws.Range("L86").Value = 132
$ws.Range("L90").Value = 216
$ws.Range("L92").Value = 62
$ws.Range("L93").Value = 103
$ws.Range("L94").Value = 250
$ws.Range("L96").Value = 227
$ws.Range("L99").Value = 360
$ws.Range("L101").Value = 20495

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L6").Value = 53
$ws.Range("L7").Value = 201

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("L2").Value = 65
$ws.Range("L7").Value = 197

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("L6").Value = 44
$ws.Range("L7").Value = 111

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 344
$ws.Range("L6").Value = 282
$ws.Range("L7").Value = 1146

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L6").Value = 109
$ws.Range("L7").Value = 269

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L6").Value = 150
$ws.Range("L7").Value = 554

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L4").Value = 37
$ws.Range("L6").Value = 142
$ws.Range("L7").Value = 317

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("L2").Value = 41
$ws.Range("L6").Value = 27
$ws.Range("L7").Value = 101

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("L3").Value = 48
$ws.Range("L7").Value = 163

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("L6").Value = 27
$ws.Range("L7").Value = 89

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L3").Value = 225
$ws.Range("L6").Value = 184
$ws.Range("L7").Value = 654

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("L2").Value = 53
$ws.Range("L7").Value = 136

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L6").Value = 77
$ws.Range("L7").Value = 268

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L6").Value = 66
$ws.Range("L7").Value = 227

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L3").Value = 182
$ws.Range("L6").Value = 122
$ws.Range("L7").Value = 522

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("L3").Value = 49
$ws.Range("L7").Value = 140

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("L3").Value = 28
$ws.Range("L7").Value = 103

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L2").Value = 230
$ws.Range("L7").Value = 658

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L6").Value = 92
$ws.Range("L7").Value = 250

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("L4").Value = 9
$ws.Range("L7").Value = 59

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L4").Value = 27
$ws.Range("L6").Value = 86
$ws.Range("L7").Value = 338

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("L2").Value = 27
$ws.Range("L7").Value = 62

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("L4").Value = 72
$ws.Range("L7").Value = 132

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L2").Value = 71
$ws.Range("L7").Value = 216

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("L3").Value = 43
$ws.Range("L7").Value = 134

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("L4").Value = 25
$ws.Range("L7").Value = 154

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 308
$ws.Range("L3").Value = 423
$ws.Range("L6").Value = 210
$ws.Range("L7").Value = 1019

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 67

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("L4").Value = 3
$ws.Range("L7").Value = 54

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L4").Value = 28
$ws.Range("L6").Value = 124
$ws.Range("L7").Value = 436
